# Auto-generated edit script applying numeric corrections to the
# Kujata_Profits workbook market-data tables (per-sheet Leve profit data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I15").Value = 2760.739
$ws.Range("K15").Value = 8282.217000000001
$ws.Range("M15").Value = -8113.217000000001
$ws.Range("H15").Value = 2760.739
$ws.Range("M69").Value = -8126
$ws.Range("N69").Value = -12218
$ws.Range("L69").Value = 10470
$ws.Range("K69").Value = 9000
$ws.Range("I69").Value = 3000
$ws.Range("H69").Value = 3392
$ws.Range("J69").Value = 3490
$ws.Range("N72").Value = -40146
$ws.Range("H72").Value = 3392
$ws.Range("K72").Value = 27000
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 3490
$ws.Range("M72").Value = -22632
$ws.Range("L72").Value = 31410
$ws.Range("L138").Value = 5479.6362
$ws.Range("J138").Value = 1826.5454
$ws.Range("H138").Value = 1584.5858
$ws.Range("N138").Value = -15759.6362
$ws.Range("I141").Value = 2316
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -1768
$ws.Range("K141").Value = 6948
$ws.Range("H141").Value = 2392
$ws.Range("J141").Value = 3000
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4424.232
$ws.Range("I32").Value = 4287.6567
$ws.Range("M32").Value = -4000.6567
$ws.Range("K32").Value = 4287.6567
$ws.Range("H44").Value = 21922
$ws.Range("K44").Value = 14044
$ws.Range("M44").Value = -13556
$ws.Range("I44").Value = 14044
$ws.Range("N45").Value = -2120.6666
$ws.Range("I45").Value = 1054
$ws.Range("J45").Value = 1366.6666
$ws.Range("M45").Value = -677
$ws.Range("L45").Value = 1366.6666
$ws.Range("K45").Value = 1054
$ws.Range("H45").Value = 1112.625
$ws.Range("N61").Value = -2757.3333
$ws.Range("K61").Value = 71429540
$ws.Range("H61").Value = 58824736
$ws.Range("J61").Value = 2333.3333
$ws.Range("I61").Value = 71429540
$ws.Range("M61").Value = -71429328
$ws.Range("L61").Value = 2333.3333
$ws.Range("H110").Value = 1350.25
$ws.Range("K110").Value = 827
$ws.Range("N110").Value = -9103
$ws.Range("J110").Value = 5013
$ws.Range("I110").Value = 827
$ws.Range("M110").Value = 1218
$ws.Range("L110").Value = 5013
$ws.Range("N132").Value = -17573
$ws.Range("H132").Value = 3093.762
$ws.Range("L132").Value = 12513
$ws.Range("K132").Value = 7988.6001
$ws.Range("I132").Value = 2662.8667
$ws.Range("M132").Value = -5458.6001
$ws.Range("J132").Value = 4171
$ws.Range("J136").Value = 2333.3333
$ws.Range("L136").Value = 6999.999899999999
$ws.Range("I136").Value = 71429540
$ws.Range("M136").Value = -214286070
$ws.Range("H136").Value = 58824736
$ws.Range("K136").Value = 214288620
$ws.Range("N136").Value = -12099.9999
$ws.Range("N139").Value = -39340
$ws.Range("H139").Value = 29060
$ws.Range("J139").Value = 29060
$ws.Range("L139").Value = 29060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J134").Value = 2228.25
$ws.Range("M134").Value = -811.0907999999999
$ws.Range("I134").Value = 1115.3636
$ws.Range("N134").Value = -11754.75
$ws.Range("H134").Value = 1412.1333
$ws.Range("K134").Value = 3346.0908
$ws.Range("L134").Value = 6684.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J135").Value = 34460
$ws.Range("L135").Value = 34460
$ws.Range("H135").Value = 34460
$ws.Range("N135").Value = -44600
$ws.Range("L141").Value = 767772.5
$ws.Range("H141").Value = 767772.5
$ws.Range("J141").Value = 767772.5
$ws.Range("N141").Value = -778132.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 638.6667
$ws.Range("N113").Value = -6411.7727
$ws.Range("L113").Value = 2071.7727
$ws.Range("M113").Value = 565.5454
$ws.Range("I113").Value = 534.8182
$ws.Range("K113").Value = 1604.4546
$ws.Range("J113").Value = 690.5909

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L2").Value = 142
$ws.Range("N2").Value = -368
$ws.Range("H2").Value = 160.14285
$ws.Range("I2").Value = 178.28572
$ws.Range("J2").Value = 142
$ws.Range("K2").Value = 178.28572
$ws.Range("M2").Value = -65.28572
$ws.Range("H53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("L97").Value = 592.4286
$ws.Range("H97").Value = 550.1053000000001
$ws.Range("M97").Value = -29.41669999999999
$ws.Range("N97").Value = -1584.4286
$ws.Range("J97").Value = 592.4286
$ws.Range("I97").Value = 525.4167
$ws.Range("K97").Value = 525.4167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N68").Value = -3898
$ws.Range("I68").Value = 2197.6667
$ws.Range("K68").Value = 2197.6667
$ws.Range("M68").Value = -1448.6667
$ws.Range("L68").Value = 2400
$ws.Range("J68").Value = 2400
$ws.Range("H68").Value = 2226.5715
$ws.Range("L71").Value = 12000
$ws.Range("I71").Value = 2197.6667
$ws.Range("M71").Value = -7244.333500000001
$ws.Range("K71").Value = 10988.3335
$ws.Range("H71").Value = 2226.5715
$ws.Range("J71").Value = 2400
$ws.Range("N71").Value = -19488
$ws.Range("M74").ClearContents()
$ws.Range("K74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("L76").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("L79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("J82").Value = 1842.1666
$ws.Range("M82").Value = -1597.1666
$ws.Range("K82").Value = 1958.1666
$ws.Range("H82").Value = 1919.5
$ws.Range("N82").Value = -2564.1666
$ws.Range("L82").Value = 1842.1666
$ws.Range("I82").Value = 1958.1666
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("I85").Value = 1958.1666
$ws.Range("H85").Value = 1919.5
$ws.Range("N85").Value = -4338.1666
$ws.Range("L85").Value = 1842.1666
$ws.Range("K85").Value = 1958.1666
$ws.Range("M85").Value = -710.1666
$ws.Range("J85").Value = 1842.1666
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("L90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("H112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
